# Applies the cryptos.xlsx data refresh described by the commit:
# "Updated cryptos list on Tue May  7 20:45:58 UTC 2024 with GitHub Actions"
# Updates Price (D) / Volume(1h) (E) figures, and fixes two mis-ordered
# coin rows (FirstDigitalUSD/PancakeSwap and Filecoin/dogwifhat).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "62.992.94"
$ws.Range("E2").Value = "  -0.78%  "

# Row 3
$ws.Range("D3").Value = "3.052.58"
$ws.Range("E3").Value = "  -0.17%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.08"
$ws.Range("E5").Value = "  -0.87%  "

# Row 6
$ws.Range("E6").Value = "  -2.00%  "

# Row 7
$ws.Range("E7").Value = "  +0.09%  "

# Row 8
$ws.Range("E8").Value = "  -1.98%  "

# Row 9
$ws.Range("D9").Value = "3.052.47"
$ws.Range("E9").Value = "  -0.78%  "

# Row 10
$ws.Range("E10").Value = "  -2.92%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.84"
$ws.Range("E11").Value = "  -0.31%  "

# Row 12
$ws.Range("E12").Value = "  -2.79%  "

# Row 13
$ws.Range("E13").Value = "  -3.19%  "

# Row 14
$ws.Range("E14").Value = "  -3.93%  "

# Row 15
$ws.Range("E15").Value = "  +1.71%  "

# Row 16
$ws.Range("D16").Value = "3.555.67"
$ws.Range("E16").Value = "  -0.24%  "

# Row 17
$ws.Range("E17").Value = "  -1.44%  "

# Row 18
$ws.Range("D18").Value = "62.991.00"
$ws.Range("E18").Value = "  -0.67%  "

# Row 19
$ws.Range("D19").Value = "3.053.54"
$ws.Range("E19").Value = "  -0.51%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "477.63"
$ws.Range("E20").Value = "  -0.23%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.26"
$ws.Range("E21").Value = "  -3.06%  "

# Row 22
$ws.Range("E22").Value = "  -1.90%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.49"
$ws.Range("E23").Value = "  -1.15%  "

# Row 24
$ws.Range("E24").Value = "  -0.81%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.86"
$ws.Range("E25").Value = "  +0.82%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.61"
$ws.Range("E26").Value = "  -2.76%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.48"
$ws.Range("E27").Value = "  +4.66%  "

# Row 28
$ws.Range("E28").Value = "  +0.27%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.36"
$ws.Range("E29").Value = "  -0.51%  "

# Row 30
$ws.Range("B30").Value = "FirstDigitalUSD"
$ws.Range("C30").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.03%  "

# Row 31
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.66"
$ws.Range("E31").Value = "  -1.15%  "

# Row 32
$ws.Range("E32").Value = "  -0.57%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.56"
$ws.Range("E33").Value = "  +1.12%  "

# Row 34
$ws.Range("E34").Value = "  -3.06%  "

# Row 35
$ws.Range("E35").Value = "  +1.02%  "

# Row 36
$ws.Range("D36").Value = "0.0₃0814"
$ws.Range("E36").Value = "  -4.16%  "

# Row 37
$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.90"
$ws.Range("E37").Value = "  -3.61%  "

# Row 38
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.23"
$ws.Range("E38").Value = "  -4.69%  "

# Row 39
$ws.Range("E39").Value = "  -1.18%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.22"
$ws.Range("E40").Value = "  -1.76%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.38"
$ws.Range("E41").Value = "  -0.52%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "431.49"
$ws.Range("E42").Value = "  -3.24%  "

# Row 43
$ws.Range("E43").Value = "  +0.62%  "

# Row 44
$ws.Range("E44").Value = "  +2.79%  "

# Row 46
$ws.Range("D46").Value = "2.829.24"
$ws.Range("E46").Value = "  +0.90%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "38.22"
$ws.Range("E47").Value = "  -4.68%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.70"
$ws.Range("E48").Value = "  -3.19%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.95"
$ws.Range("E50").Value = "  -1.24%  "

# Row 51
$ws.Range("E51").Value = "  -1.61%  "
